$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contract")
$ws.Range("B27").Value = "-"
$ws.Range("B28").Value = "-"
$ws.Range("B29").Value = "-"
$ws.Range("B30").Value = "-"
$ws.Range("B31").Value = "-"
$ws.Range("B32").Value = "-"
